# "Some better resource management"
# Add a new "Oxygen Producing Unit" resource row to the Equipment/Upgrades
# table on Sheet2, just above the existing "Passive H3 Unit" row, and move
# the selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert a new row before the current row 4 ("Passive H3 Unit"); this shifts
# every row below it down by one and keeps all relative formulas intact.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new resource's data.
$ws.Range("A4").Value = "Oxygen Producing Unit"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 10

# Leave the selection where the author left it.
[void]$ws.Range("C3").Select()
